# "Went to 2nd order" - append 7 new data rows (31-37) to the weight-tracking
# sheet, continuing the existing A (oz), B (date), C (=B-$B$2 DeltaDays) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows: row, A (oz), B (date serial)
$newRows = @(
    @(31, 5.6, 44561),
    @(32, 6.1, 44562),
    @(33, 7.5, 44563),
    @(34, 8.9, 44565),
    @(35, 8.4, 44566),
    @(36, 4.9, 44570),
    @(37, 4.5, 44572)
)

foreach ($r in $newRows) {
    $row = $r[0]
    $oz = $r[1]
    $dateSerial = $r[2]

    $ws.Cells.Item($row, 1).Value = $oz

    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = $dateSerial
    $bCell.NumberFormat = "d-mmm"

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Formula = "=B$row-`$B`$2"
    # Setting the formula picks up the date-like number format from the
    # referenced date cells; the source column is plain/General, so reset it.
    $cCell.Style = "Normal"
}

# Match the saved selection/view state.
$null = $ws.Range("B38").Select()
